# Applies the crypto price/volume refresh described in the commit.
# D-column (Price) values are forced to Text so numeric-looking
# strings (e.g. "216.41") are not auto-coerced into numbers by Excel,
# matching the original inlineStr cell type. The style is reset back
# to "Normal" afterwards so no stray number-format style is left on
# the cell (only the cells value type needs to become text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "26.021.78"
$ws.Range("E2").Value = "  -2.28%  "

# Row 3
Set-TextValue "D3" "1.665.04"
$ws.Range("E3").Value = "  -1.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
Set-TextValue "D5" "216.41"
$ws.Range("E5").Value = "  -1.04%  "

# Row 6
Set-TextValue "D6" "0.5089"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("E7").Value = "  +0.23%  "

# Row 8
$ws.Range("E8").Value = "  -0.55%  "

# Row 9
Set-TextValue "D9" "0.06391"
$ws.Range("E9").Value = "  +2.27%  "

# Row 10
Set-TextValue "D10" "21.69"
$ws.Range("E10").Value = "  -1.99%  "

# Row 11
Set-TextValue "D11" "0.07401"
$ws.Range("E11").Value = "  +1.32%  "

# Row 12
Set-TextValue "D12" "1.670.48"
$ws.Range("E12").Value = "  -1.48%  "

# Row 13
Set-TextValue "D13" "4.498"
$ws.Range("E13").Value = "  +0.41%  "

# Row 14
Set-TextValue "D14" "0.5805"
$ws.Range("E14").Value = "  -0.26%  "

# Row 15
Set-TextValue "D15" "0.000008505"
$ws.Range("E15").Value = "  +2.22%  "

# Row 16
Set-TextValue "D16" "64.17"
$ws.Range("E16").Value = "  -1.82%  "

# Row 17
Set-TextValue "D17" "26.093.84"
$ws.Range("E17").Value = "  -2.13%  "

# Row 18
$ws.Range("E18").Value = "  -2.87%  "

# Row 19
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
Set-TextValue "D20" "10.69"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
Set-TextValue "D21" "188.59"
$ws.Range("E21").Value = "  +1.22%  "

# Row 22
Set-TextValue "D22" "6.197"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
Set-TextValue "D24" "145.72"
$ws.Range("E24").Value = "  +0.51%  "

# Row 25
Set-TextValue "D25" "7.601"
$ws.Range("E25").Value = "  +0.31%  "

# Row 26
Set-TextValue "D26" "0.1187"
$ws.Range("E26").Value = "  +3.91%  "

# Row 27
Set-TextValue "D27" "15.59"
$ws.Range("E27").Value = "  -0.43%  "

# Row 28
Set-TextValue "D28" "0.06651"
$ws.Range("E28").Value = "  +16.45%  "

# Row 29
Set-TextValue "D29" "1.305"
$ws.Range("E29").Value = "  +0.69%  "

# Row 30
Set-TextValue "D30" "1.313"
$ws.Range("E30").Value = "  -1.25%  "

# Row 31
Set-TextValue "D31" "3.521"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
Set-TextValue "D33" "1.625"
$ws.Range("E33").Value = "  -1.93%  "

# Row 34
Set-TextValue "D34" "1.018"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
Set-TextValue "D35" "0.6051"
$ws.Range("E35").Value = "  +1.39%  "

# Row 36
Set-TextValue "D36" "2.367"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
Set-TextValue "D38" "6.199"

# Row 39
Set-TextValue "D39" "0.01609"
$ws.Range("E39").Value = "  +0.64%  "

# Row 40
Set-TextValue "D40" "1.073.87"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
Set-TextValue "D41" "0.8591"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("E42").Value = "  +0.66%  "

# Row 43
Set-TextValue "D43" "100.40"
$ws.Range("E43").Value = "  +1.98%  "

# Row 44
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D44" "0.00000000116"
$ws.Range("E44").Value = "  +8.95%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D45" "1.812.74"
$ws.Range("E45").Value = "  -2.25%  "

# Row 46
Set-TextValue "D46" "56.19"
$ws.Range("E46").Value = "  -0.89%  "

# Row 47
Set-TextValue "D47" "1.008"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48
Set-TextValue "D48" "7.991"
$ws.Range("E48").Value = "  -1.54%  "

# Row 49
Set-TextValue "D49" "0.05208"
$ws.Range("E49").Value = "  -0.58%  "

# Row 50
Set-TextValue "D50" "0.4292"
$ws.Range("E50").Value = "  -0.50%  "

# Row 51
Set-TextValue "D51" "5.944"
$ws.Range("E51").Value = "  +3.87%  "
